$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (IDPT) updates
$ws.Range("B2").Value = 0.95649861250216
$ws.Range("C2").Value = 0.5174819885839491
$ws.Range("G2").Value = 0.5248562840851071
$ws.Range("H2").Value = 0.4794191554825576
$ws.Range("I2").Value = 0.7108564169980006
$ws.Range("J2").Value = 0.4368212139013902
$ws.Range("K2").Value = 0.8343438251113983

# Row 3 label update: SPCT -> GDPT
$ws.Range("A3").Value = "GDPT"

# Row 3 (GDPT) numeric updates
$ws.Range("C3").Value = 1.615307186802609
$ws.Range("G3").Value = 0.3383571169547176
$ws.Range("H3").Value = 0.4136789058781119
$ws.Range("I3").Value = 0.5344303282584365
$ws.Range("J3").Value = 1.055971845676007
$ws.Range("K3").Value = 1.183508476785364
